$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Widen a few week-marker columns on the Gantt chart so the dates/labels are legible.
$ws.Columns.Item(22).ColumnWidth = 13.666666666666666   # V  -> ~14.42578125
$ws.Columns.Item(29).ColumnWidth = 11.666666666666666   # AC -> ~12.5703125
$ws.Columns.Item(35).ColumnWidth = 10.666666666666666   # AI -> ~11.42578125
$ws.Columns.Item(42).ColumnWidth = 12.0                 # AP -> ~12.85546875

# Zoom out and move the selection, as in the saved view state.
$excel.ActiveWindow.Zoom = 60
$ws.Range("AH2").Select()
